$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.975.89"
$ws.Range("E2").Value = "  +3.46%  "
$ws.Range("D3").Value = "2.587.57"
$ws.Range("E3").Value = "  +2.81%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'323.08"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").Value = "'109.62"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").Value = "'0.531"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "'0.560"
$ws.Range("E9").Value = "  +2.68%  "
$ws.Range("D10").Value = "'40.65"
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("D11").Value = "'20.53"
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("D12").Value = "'0.0821"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "'7.29"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("D15").Value = "2.966.30"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "2.562.77"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").Value = "'0.866"
$ws.Range("E17").Value = "  +2.27%  "
$ws.Range("D18").Value = "49.787.46"
$ws.Range("E18").Value = "  +3.41%  "
$ws.Range("D19").Value = "'3.09"
$ws.Range("E19").Value = "  +12.55%  "
$ws.Range("D20").Value = "'13.41"
$ws.Range("E20").Value = "  +1.91%  "
$ws.Range("D21").Value = "'6.74"
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("D22").Value = "0.0₃0951"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "'284.62"
$ws.Range("E23").Value = "  +2.48%  "
$ws.Range("D24").Value = "'72.75"
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("D25").Value = "'2.53"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("D26").Value = "'26.65"
$ws.Range("E26").Value = "  +2.66%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  +4.81%  "
$ws.Range("E29").Value = "  -7.13%  "
$ws.Range("D30").Value = "'9.97"
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("D31").Value = "'36.19"
$ws.Range("E31").Value = "  +1.96%  "
$ws.Range("D32").Value = "'49.54"
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("D33").Value = "'19.76"
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("D34").Value = "'5.43"
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").Value = "'0.0791"
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("D37").Value = "'2.05"
$ws.Range("E37").Value = "  +4.16%  "
$ws.Range("D38").Value = "'4.78"
$ws.Range("E38").Value = "  +2.97%  "
$ws.Range("D39").Value = "'3.06"
$ws.Range("E39").Value = "  +3.08%  "
$ws.Range("D40").Value = "'123.61"
$ws.Range("E40").Value = "  +1.12%  "
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("D42").Value = "'22.30"
$ws.Range("E42").Value = "  +3.49%  "
$ws.Range("D43").Value = "'2.22"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").Value = "'0.0314"
$ws.Range("E44").Value = "  +2.70%  "
$ws.Range("E45").Value = "  +5.67%  "
$ws.Range("D46").Value = "2.028.84"
$ws.Range("E46").Value = "  +1.28%  "
$ws.Range("E47").Value = "  +9.43%  "
$ws.Range("D48").Value = "'2.17"
$ws.Range("E48").Value = "  +9.13%  "
$ws.Range("D49").Value = "'9.13"
$ws.Range("E49").Value = "  +0.89%  "
$ws.Range("D50").Value = "'5.38"
$ws.Range("E50").Value = "  +2.44%  "
$ws.Range("D51").Value = "'81.95"
$ws.Range("E51").Value = "  +2.07%  "
